$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (ownTeam, oppTeam) before the existing "batsman"
# column, pushing batsman/totalRuns/totalBalls/total4s/total6s/sr from
# D:I to F:K.
$ws.Range("D1:E1").EntireColumn.Insert()

# --- Header row ---
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# --- Row 2 (existing match) gets the two new team columns filled in ---
$ws.Range("D2").Value = "Kings XI Punjab"
$ws.Range("E2").Value = "Chennai Super Kings"

# Re-enter the numeric-looking stat columns as text so they keep being
# stored as text (matches the source data, which stores every value as
# text, never as a real number).
$ws.Range("G2:K2").NumberFormat = "@"
$ws.Range("G2").Value = "14"
$ws.Range("H2").Value = "9"
$ws.Range("I2").Value = "2"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "155.55"

# --- New row 3 ---
$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " October 01 2020"
$ws.Range("C3").Value = "Mumbai won by 48 runs"
$ws.Range("D3").Value = "Kings XI Punjab"
$ws.Range("E3").Value = "Mumbai Indians"
$ws.Range("F3").Value = "Sarfaraz Khan" + [char]0x00A0
$ws.Range("G3:K3").NumberFormat = "@"
$ws.Range("G3").Value = "7"
$ws.Range("H3").Value = "8"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "87.50"

# --- New row 4 ---
$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " September 20 2020"
$ws.Range("C4").Value = "Match tied (Capitals won the one-over eliminator)"
$ws.Range("D4").Value = "Kings XI Punjab"
$ws.Range("E4").Value = "Delhi Capitals"
$ws.Range("F4").Value = "Sarfaraz Khan" + [char]0x00A0
$ws.Range("G4:K4").NumberFormat = "@"
$ws.Range("G4").Value = "12"
$ws.Range("H4").Value = "12"
$ws.Range("I4").Value = "2"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "100.00"
